$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values for rows 2-5 (columns B..F)
$ws.Range("B2").Value2 = -0.2902542872142831
$ws.Range("C2").Value2 = 0.4474900817361138
$ws.Range("D2").Value2 = 0.3638477001211933
$ws.Range("E2").Value2 = 0.6031978946591187
$ws.Range("F2").Value2 = 0.5487329520217746

$ws.Range("B3").Value2 = -0.08670396990747065
$ws.Range("C3").Value2 = 0.4011283099121256
$ws.Range("D3").Value2 = 0.2437355075710951
$ws.Range("E3").Value2 = 0.4936957641818442
$ws.Range("F3").Value2 = 0.5123127616487798
$ws.Range("G3").Value2 = 10

$ws.Range("B4").Value2 = -0.2842180920556761
$ws.Range("C4").Value2 = 0.4543513181817999
$ws.Range("D4").Value2 = 0.2495528846253024
$ws.Range("E4").Value2 = 0.4995526845341764
$ws.Range("F4").Value2 = 0.4500306133234053
$ws.Range("G4").Value2 = 6

$ws.Range("B5").Value2 = -0.2239409044200031
$ws.Range("C5").Value2 = 0.5803388081800324
$ws.Range("D5").Value2 = 0.3869426609522694
$ws.Range("E5").Value2 = 0.6220471533189983
$ws.Range("F5").Value2 = 0.8207230132996398
$ws.Range("G5").Value2 = 2

# Remove rows 6-9 (Q4..Q7), which also prunes the now-unused shared strings
$ws.Rows("6:9").Delete()
